$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 95290
$ws.Range("B2").Value = "Laura Barbosa"
$ws.Range("C2").Value = "Recursos Humanos"
$ws.Range("D2").Value = "Consulta médica"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 45097
$ws.Range("G2").Value = 11434.27

$ws.Range("A3").Value = 6340
$ws.Range("B3").Value = "Danilo Jesus"
$ws.Range("C3").Value = "P&D"
$ws.Range("E3").Value = 7
$ws.Range("F3").Value = 45093
$ws.Range("G3").Value = 11539.98

$ws.Range("A4").Value = 27962
$ws.Range("B4").Value = "Isadora da Mata"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 45094
$ws.Range("G4").Value = 11486.89

$ws.Range("A5").Value = 75376
$ws.Range("B5").Value = "Elisa Carvalho"
$ws.Range("C5").Value = "Atendimento ao Cliente"
$ws.Range("E5").Value = 5
$ws.Range("F5").Value = 45106
$ws.Range("G5").Value = 10313.17

$ws.Range("A6").Value = 9985
$ws.Range("B6").Value = "Luiz Gustavo Rodrigues"
$ws.Range("C6").Value = "Marketing"
$ws.Range("D6").Value = "Viagem de negócios"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 45093
$ws.Range("G6").Value = 10797.29

$ws.Range("A7").Value = 44274
$ws.Range("B7").Value = "Dra. Júlia Ribeiro"
$ws.Range("C7").Value = "Vendas"
$ws.Range("D7").Value = "Doença"
$ws.Range("E7").Value = 6
$ws.Range("F7").Value = 45086
$ws.Range("G7").Value = 5322.27

$ws.Range("A8").Value = 52539
$ws.Range("B8").Value = "Sarah Oliveira"
$ws.Range("C8").Value = "Atendimento ao Cliente"
$ws.Range("D8").Value = "Doença"
$ws.Range("E8").Value = 6
$ws.Range("F8").Value = 45080
$ws.Range("G8").Value = 9475.85

$ws.Range("A9").Value = 60648
$ws.Range("B9").Value = "Srta. Sophia da Luz"
$ws.Range("C9").Value = "Vendas"
$ws.Range("D9").Value = "Outros"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 45094
$ws.Range("G9").Value = 10986.13

$ws.Range("A10").Value = 26197
$ws.Range("B10").Value = "Lorenzo Teixeira"
$ws.Range("C10").Value = "Engenharia"
$ws.Range("D10").Value = "Doença"
$ws.Range("F10").Value = 45096
$ws.Range("G10").Value = 11623.72

$ws.Range("A11").Value = 67535
$ws.Range("B11").Value = "Leandro Costela"
$ws.Range("C11").Value = "Jurídico"
$ws.Range("D11").Value = "Viagem de negócios"
$ws.Range("E11").Value = 4
$ws.Range("F11").Value = 45094
$ws.Range("G11").Value = 4593.65
